$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.807.52"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "3.339.70"
$ws.Range("E3").Value = "  +2.68%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.96"
$ws.Range("E5").Value = "  +2.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.26"
$ws.Range("E6").Value = "  -3.23%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "3.344.81"
$ws.Range("E8").Value = "  +3.17%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.605"
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.10"
$ws.Range("E11").Value = "  -6.38%  "
$ws.Range("E12").Value = "  +3.54%  "
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.06"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "3.871.50"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("D16").Value = "3.337.27"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.51"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "63.739.94"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.14"
$ws.Range("E20").Value = "  +3.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.961"
$ws.Range("E21").Value = "  +2.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "372.58"
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.46"
$ws.Range("E23").Value = "  +3.18%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.07"
$ws.Range("E24").Value = "  +8.15%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.42"
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.70"
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.18"
$ws.Range("E27").Value = "  +3.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.69"
$ws.Range("E28").Value = "  +3.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.26"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.21"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.72"
$ws.Range("E31").Value = "  +2.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "635.84"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.39"
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.16"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.51"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.26"
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("D40").Value = "0.0₃0718"
$ws.Range("E40").Value = "  +11.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.997"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("E42").Value = "  +8.46%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.125"
$ws.Range("E43").Value = "  +1.64%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.939.36"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("E45").Value = "  +9.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.67"
$ws.Range("E46").Value = "  +2.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0394"
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.59"
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.05"
$ws.Range("E49").Value = "  +4.18%  "
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.39"
$ws.Range("E51").Value = "  +4.94%  "
